$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '45.875.02'
Set-TextCell 'E2' '  -2.47%  '
Set-TextCell 'D3' '2.284.65'
Set-TextCell 'E3' '  -2.44%  '
Set-TextCell 'D4' '0.998'
Set-TextCell 'E4' '  -0.11%  '
Set-TextCell 'D5' '297.03'
Set-TextCell 'E5' '  -2.93%  '
Set-TextCell 'D6' '98.04'
Set-TextCell 'E6' '  +0.44%  '
Set-TextCell 'D7' '0.567'
Set-TextCell 'E7' '  -1.76%  '
Set-TextCell 'D8' '0.999'
Set-TextCell 'D9' '0.505'
Set-TextCell 'E9' '  -6.03%  '
Set-TextCell 'D10' '34.46'
Set-TextCell 'E10' '  -4.09%  '
Set-TextCell 'D11' '0.0775'
Set-TextCell 'E11' '  -4.22%  '
Set-TextCell 'D12' '6.97'
Set-TextCell 'E12' '  -6.26%  '
Set-TextCell 'E13' '  -1.97%  '
Set-TextCell 'D14' '2.620.90'
Set-TextCell 'E14' '  -2.80%  '
Set-TextCell 'D15' '2.306.61'
Set-TextCell 'E15' '  -1.47%  '
Set-TextCell 'D16' '13.50'
Set-TextCell 'E16' '  -4.94%  '
Set-TextCell 'D17' '0.791'
Set-TextCell 'E17' '  -4.86%  '
Set-TextCell 'D18' '45.802.61'
Set-TextCell 'E18' '  -2.17%  '
Set-TextCell 'D19' '12.39'
Set-TextCell 'E19' '  -8.43%  '
Set-TextCell 'D20' '0.0₃0953'
Set-TextCell 'E20' '  +0.31%  '
Set-TextCell 'D21' '5.79'
Set-TextCell 'E21' '  -6.48%  '
Set-TextCell 'D22' '65.10'
Set-TextCell 'E22' '  -3.36%  '
Set-TextCell 'D23' '242.48'
Set-TextCell 'E23' '  -2.12%  '
Set-TextCell 'D24' '2.77'
Set-TextCell 'E24' '  -7.12%  '
Set-TextCell 'E25' '  +0.37%  '
Set-TextCell 'D26' '1.85'
Set-TextCell 'E26' '  -7.65%  '
Set-TextCell 'D27' '40.15'
Set-TextCell 'E27' '  -6.37%  '
Set-TextCell 'E28' '  -4.29%  '
Set-TextCell 'D29' '9.50'
Set-TextCell 'E29' '  -4.26%  '
Set-TextCell 'D30' '19.94'
Set-TextCell 'E30' '  -1.15%  '
Set-TextCell 'D31' '2.80'
Set-TextCell 'E31' '  +7.10%  '
Set-TextCell 'B32' 'LidoDAOToken'
Set-TextCell 'C32' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 'D32' '3.34'
Set-TextCell 'E32' '  +4.94%  '
Set-TextCell 'B33' 'Monero'
Set-TextCell 'C33' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D33' '143.41'
Set-TextCell 'E33' '  -3.59%  '
Set-TextCell 'B34' 'Filecoin'
Set-TextCell 'C34' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D34' '5.28'
Set-TextCell 'E34' '  -8.72%  '
Set-TextCell 'D35' '0.0763'
Set-TextCell 'E35' '  -6.92%  '
Set-TextCell 'D36' '0.110'
Set-TextCell 'E36' '  -2.69%  '
Set-TextCell 'E37' '  -3.88%  '
Set-TextCell 'D38' '15.28'
Set-TextCell 'E38' '  +9.48%  '
Set-TextCell 'E39' '  -9.07%  '
Set-TextCell 'D40' '3.78'
Set-TextCell 'E40' '  -5.92%  '
Set-TextCell 'D41' '0.0293'
Set-TextCell 'E41' '  -7.12%  '
Set-TextCell 'D42' '3.08'
Set-TextCell 'E42' '  -9.62%  '
Set-TextCell 'D43' '0.998'
Set-TextCell 'E43' '  -0.05%  '
Set-TextCell 'D44' '93.47'
Set-TextCell 'E44' '  +10.01%  '
Set-TextCell 'D45' '1.808.47'
Set-TextCell 'E45' '  -2.04%  '
Set-TextCell 'D46' '1.84'
Set-TextCell 'E46' '  -7.95%  '
Set-TextCell 'B47' 'Algorand'
Set-TextCell 'C47' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D47' '0.182'
Set-TextCell 'E47' '  -7.21%  '
Set-TextCell 'B48' 'ordi'
Set-TextCell 'C48' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextCell 'D48' '69.08'
Set-TextCell 'E48' '  -8.64%  '
Set-TextCell 'B49' 'RocketPoolETH'
Set-TextCell 'C49' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 'D49' '2.504.94'
Set-TextCell 'E49' '  -2.65%  '
Set-TextCell 'B50' 'THORChain'
Set-TextCell 'C50' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell 'D50' '4.67'
Set-TextCell 'E50' '  -4.91%  '
Set-TextCell 'B51' 'Aave'
Set-TextCell 'C51' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 'D51' '93.97'
Set-TextCell 'E51' '  -5.31%  '
